$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "AL-SQP"

$ws.Range("B2").Value = 493.2489079521391
$ws.Range("C2").Value = 398.0125759009749
$ws.Range("D2").Value = 489.6521082785899
$ws.Range("E2").Value = 411.9903044170762
$ws.Range("F2").Value = 498.5425510850154

$ws.Range("B3").Value = 490.92598839141
$ws.Range("C3").Value = 420.1427919264566
$ws.Range("D3").Value = 493.364789575227
$ws.Range("E3").Value = 403.6789918291935
$ws.Range("F3").Value = 499.5835623833169

$ws.Range("B4").Value = 486.5966483399661
$ws.Range("C4").Value = 359.1921061179575
$ws.Range("D4").Value = 491.2368821989236
$ws.Range("E4").Value = 283.4839307807926
$ws.Range("F4").Value = 495.0859409474938

$ws.Range("B5").Value = 466.6824725294445
$ws.Range("C5").Value = 398.9847710574171
$ws.Range("D5").Value = 477.0345653802399
$ws.Range("E5").Value = 389.2410313257318
$ws.Range("F5").Value = 479.2495412993596

$ws.Range("B6").Value = 465.4623654255223
$ws.Range("C6").Value = 364.838061883233
$ws.Range("D6").Value = 480.2212943784397
$ws.Range("E6").Value = 334.1144620973773
$ws.Range("F6").Value = 483.4524913894659

$ws.Range("B7").Value = 452.6495762701873
$ws.Range("C7").Value = 364.8549915166554
$ws.Range("D7").Value = 471.4130381358773
$ws.Range("E7").Value = 391.4503704283915
$ws.Range("F7").Value = 475.007333934248

$ws.Range("B8").Value = 472.3969636003201
$ws.Range("C8").Value = 403.881023021193
$ws.Range("D8").Value = 494.9636192133948
$ws.Range("E8").Value = 390.1819007836912
$ws.Range("F8").Value = 497.9661947390633
